$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: target row, D, M, N, O, P, S (derived from row-shuffle in source diff)
$rowUpdates = @(
    @(2, 44753, 160, 2300, 2300, 2300, 2300),
    @(3, 45041, 80, 3500, 3500, 3500, 3500),
    @(4, 45068, 50, 3250, 3250, 3250, 3250),
    @(5, 45104, 50, 2600, 2600, 2600, 2600),
    @(6, 44418, 40, 1200, 1200, 1200, 1200),
    @(7, 44830, 50, 2500, 2500, 2500, 2500),
    @(8, 44760, 80, 2300, 2300, 2300, 2300),
    @(9, 44343, 60, 1300, 1300, 1300, 1300),
    @(10, 44476, 80, 1200, 1200, 1200, 1200),
    @(11, 44432, 30, 1300, 1300, 1300, 1300),
    @(12, 45042, 25, 3500, 3500, 3500, 3500),
    @(13, 45076, 100, 2600, 2600, 2600, 2600),
    @(14, 45086, 30, 2600, 2600, 2600, 2600),
    @(15, 44357, 35, 1000, 1000, 1000, 1000),
    @(16, 45093, 90, 2600, 2600, 2600, 2600),
    @(17, 45092, 120, 2600, 2600, 2600, 2600),
    @(18, 45044, 150, 3500, 3500, 3500, 3500),
    @(19, 45054, 25, 2500, 2500, 2500, 2500),
    @(20, 44762, 50, 2300, 2300, 2300, 2300),
    @(21, 44405, 50, 1200, 1200, 1200, 1200),
    @(22, 45113, 90, 2600, 2600, 2600, 2600),
    @(23, 44748, 300, 2300, 2300, 2300, 2300),
    @(24, 44431, 100, 1300, 1300, 1300, 1300),
    @(25, 44435, 130, 1300, 1300, 1300, 1300),
    @(26, 44473, 120, 1200, 1200, 1200, 1200),
    @(27, 44812, 50, 2500, 2500, 2500, 2500),
    @(28, 44424, 50, 1200, 1200, 1200, 1200),
    @(29, 45097, 90, 2600, 2600, 2600, 2600),
    @(30, 45090, 50, 2600, 2600, 2600, 2600),
    @(31, 45085, 40, 2600, 2600, 2600, 2600),
    @(32, 44438, 60, 1200, 1200, 1200, 1200),
    @(33, 45055, 25, 2800, 2800, 2800, 2800),
    @(34, 45062, 60, 3200, 3200, 3200, 3200),
    @(35, 45079, 30, 2600, 2600, 2600, 2600),
    @(36, 45111, 50, 2600, 2600, 2600, 2600),
    @(37, 44417, 80, 1200, 1200, 1200, 1200),
    @(38, 45106, 120, 2600, 2600, 2600, 2600),
    @(39, 45075, 240, 3200, 3200, 3200, 3200),
    @(40, 44811, 60, 2500, 2500, 2500, 2500),
    @(41, 45112, 50, 2600, 2600, 2600, 2600),
    @(42, 44749, 120, 2300, 2300, 2300, 2300),
    @(43, 44763, 50, 2300, 2300, 2300, 2300),
    @(44, 45099, 200, 2600, 2600, 2600, 2600)
)

foreach ($u in $rowUpdates) {
    $r = $u[0]
    $ws.Cells.Item($r, 4).Value = $u[1]   # D - Fecha
    $ws.Cells.Item($r, 13).Value = $u[2]  # M - Volumen
    $ws.Cells.Item($r, 14).Value = $u[3]  # N - Precio minimo
    $ws.Cells.Item($r, 15).Value = $u[4]  # O - Precio maximo
    $ws.Cells.Item($r, 16).Value = $u[5]  # P - Precio promedio ponderado
    $ws.Cells.Item($r, 19).Value = $u[6]  # S - Precio $/Kg
}
